# fix: Implement code review suggestions
# Adds an "e.g. ..." example line as a new paragraph under each of the
# three DevSecOps pipeline scanning stage labels.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Container Image Scanning" -> + "e.g. Amazon ECR"
$shp1 = $s.Shapes.Item(8)
[void]$shp1.TextFrame.TextRange.InsertAfter("`re.g. Amazon ECR")

# "Infrastructure as Code Scanning" -> + "e.g. Checkov"
$shp2 = $s.Shapes.Item(10)
[void]$shp2.TextFrame.TextRange.InsertAfter("`re.g. ")
[void]$shp2.TextFrame.TextRange.InsertAfter("Checkov")

# "Post-deployment compliance scanning" -> + "e.g. Amazon Inspector"
$shp3 = $s.Shapes.Item(14)
[void]$shp3.TextFrame.TextRange.InsertAfter("`re.g. ")
[void]$shp3.TextFrame.TextRange.InsertAfter("Amazon Inspector")
